$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.381.64'
$ws.Range('E2').Value = '  +1.74%  '
$ws.Range('D3').Value = '1.840.34'
$ws.Range('E3').Value = '  +1.33%  '
$ws.Range('D4').Value = '1.015'
$ws.Range('E4').Value = '  +1.35%  '
$ws.Range('D5').Value = '314.90'
$ws.Range('E5').Value = '  +1.84%  '
$ws.Range('E6').Value = '  +1.21%  '
$ws.Range('D7').Value = '0.4744'
$ws.Range('E7').Value = '  +1.59%  '
$ws.Range('D8').Value = '0.3702'
$ws.Range('E8').Value = '  +0.35%  '
$ws.Range('D9').Value = '0.07470'
$ws.Range('E9').Value = '  +1.45%  '
$ws.Range('D10').Value = '0.8856'
$ws.Range('E10').Value = '  +1.60%  '
$ws.Range('D11').Value = '20.50'
$ws.Range('E11').Value = '  +0.49%  '
$ws.Range('D12').Value = '1.851.81'
$ws.Range('E12').Value = '  +2.77%  '
$ws.Range('D13').Value = '0.07377'
$ws.Range('E13').Value = '  +4.30%  '
$ws.Range('D14').Value = '5.485'
$ws.Range('E14').Value = '  +1.96%  '
$ws.Range('D15').Value = '93.28'
$ws.Range('D16').Value = '6.585'
$ws.Range('E16').Value = '  +1.07%  '
$ws.Range('D17').Value = '1.014'
$ws.Range('E17').Value = '  +1.33%  '
$ws.Range('D18').Value = '0.000008853'
$ws.Range('E18').Value = '  +1.75%  '
$ws.Range('D19').Value = '1.014'
$ws.Range('E19').Value = '  +1.33%  '
$ws.Range('E20').Value = '  +0.85%  '
$ws.Range('D21').Value = '27.403.85'
$ws.Range('E21').Value = '  +1.76%  '
$ws.Range('D22').Value = '5.355'
$ws.Range('E22').Value = '  +0.61%  '
$ws.Range('D23').Value = '10.72'
$ws.Range('E23').Value = '  +0.94%  '
$ws.Range('D24').Value = '2.069.79'
$ws.Range('E24').Value = '  +2.33%  '
$ws.Range('D25').Value = '1.916'
$ws.Range('E25').Value = '  +1.12%  '
$ws.Range('D26').Value = '152.41'
$ws.Range('E26').Value = '  +1.29%  '
$ws.Range('D27').Value = '18.64'
$ws.Range('E27').Value = '  +1.59%  '
$ws.Range('D28').Value = '2.169'
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('D29').Value = '5.264'
$ws.Range('E29').Value = '  -1.44%  '
$ws.Range('E30').Value = '  +1.76%  '
$ws.Range('D31').Value = '0.08971'
$ws.Range('E31').Value = '  +0.55%  '
$ws.Range('D32').Value = '0.7603'
$ws.Range('E32').Value = '  -1.19%  '
$ws.Range('D33').Value = '1.181'
$ws.Range('E33').Value = '  +1.44%  '
$ws.Range('D34').Value = '4.562'
$ws.Range('E34').Value = '  +1.23%  '
$ws.Range('D35').Value = '2.940'
$ws.Range('D36').Value = '1.014'
$ws.Range('E36').Value = '  +1.38%  '
$ws.Range('D37').Value = '1.106'
$ws.Range('E37').Value = '  +1.93%  '
$ws.Range('D38').Value = '0.05380'
$ws.Range('E38').Value = '  +1.71%  '
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('E40').Value = '  +2.41%  '
$ws.Range('D41').Value = '7.283'
$ws.Range('E41').Value = '  +0.29%  '
$ws.Range('D42').Value = '0.5355'
$ws.Range('E42').Value = '  +0.52%  '
$ws.Range('D43').Value = '2.390'
$ws.Range('E43').Value = '  +1.55%  '
$ws.Range('E44').Value = '  +0.33%  '
$ws.Range('D45').Value = '8.550'
$ws.Range('E45').Value = '  +1.27%  '
$ws.Range('D46').Value = '0.4978'
$ws.Range('E46').Value = '  +0.93%  '
$ws.Range('D47').Value = '10.53'
$ws.Range('E47').Value = '  +0.41%  '
$ws.Range('D49').Value = '105.37'
$ws.Range('E49').Value = '  +1.46%  '
$ws.Range('D50').Value = '1.681'
$ws.Range('E50').Value = '  +0.48%  '
$ws.Range('D51').Value = '0.06325'
$ws.Range('E51').Value = '  +0.46%  '
